$wb = $excel.ActiveWorkbook

# Add the new worksheet right after the last existing sheet so it lands
# at the end of the tab order, then name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "test_data1"

# Write the data in a loop.
$values = @("keys", "values")
for ($i = 0; $i -lt $values.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $values[$i]
}
